$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1. Section 3 "Adopted Technologies" paragraph text
# ----------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "List the technologies used and why we chose them.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Android Studio: We have chosen to use this application since our group is familiar with it. We all feel comfortable coding in Java and installation of Android Studio is simple. Our entire group can access it with ease and the version control it implements is simple to use with GitHub. Finally, our program that we are creating is best used on a portable device. With Android Studio, creating our product for distribution to Android devices is exactly what the program was intended to do.",
    2)

# ----------------------------------------------------------------------
# 2. Section 4 "Licensing" paragraph text
# ----------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Include the license we used for our source code.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Apache License 2.0",
    2)

# ----------------------------------------------------------------------
# 3. Section 6 "Learning/Training" paragraph text (was the placeholder
#    "Describe strategies used to learn our Adopted Technologies.")
# ----------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Describe strategies used to learn our Adopted Technologies.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Android Studio has plenty of documentation which we can refer to anytime at (https://developer.android.com/guide/index.html). We also have other resources such as stackoverflow to answer questions and YouTube tutorial videos. Most of our group has also taken the mobile application course here at NAU and have used Android Studio for several of the courses projects.",
    2)

# ----------------------------------------------------------------------
# 4. Drop the stale <w:lastRenderedPageBreak/> marker that used to sit
#    in front of "Describe what we learned..." (section 7). Touching the
#    run's text via Find/Replace (even with identical text) clears it.
# ----------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Describe what we learned from this first release and what we plan on changing for the second release.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Describe what we learned from this first release and what we plan on changing for the second release.",
    2)

# ----------------------------------------------------------------------
# 5. Append "Wrote sections three and six." to Justin Shaner's line at
#    the very end of the document, and move the _GoBack bookmark there
#    (it always tracks the most recent edit location).
#
#    "Justin Shaner: " also occurs earlier (the italic example line), so
#    anchor the search after "Jacob Lemon:" to land on the real, final
#    occurrence near the end of the document.
# ----------------------------------------------------------------------
$anchor = $d.Content
$null = $anchor.Find.Execute("Jacob Lemon:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rest = $d.Range($anchor.End, $d.Content.End)
$null = $rest.Find.Execute("Justin Shaner: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insertStart = $rest.End
$ip = $d.Range($insertStart, $insertStart)
$ip.InsertAfter("Wrote sections three and six.")

# The new run lands at the very end of the body where there is no
# following character to inherit formatting from, so set it explicitly
# to match the surrounding "Book Antiqua" / 12pt (sz 24 half-points) text.
$newRun = $d.Range($insertStart, $insertStart + 30)
$newRun.Font.Name = "Book Antiqua"
$newRun.Font.Size = 12

$goBackPoint = $d.Range($d.Content.End, $d.Content.End)
$null = $d.Bookmarks.Add("_GoBack", $goBackPoint)
